$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Category"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 5)
    $val = $cell.Text
    if ($val -ne $null -and $val -ne "") {
        $cell.Value = "category " + $val
    }
}
